$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.145.57'
$ws.Range("E2").Value = '  -3.41%  '
$ws.Range("D3").Value = '3.151.72'
$ws.Range("E3").Value = '  -2.92%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.82'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.22'
$ws.Range("E6").Value = '  -6.10%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.144.62'
$ws.Range("E8").Value = '  -3.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  -3.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  -6.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  -7.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  -4.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  -7.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.49'
$ws.Range("E14").Value = '  -8.98%  '
$ws.Range("D15").Value = '3.662.38'
$ws.Range("E15").Value = '  -3.10%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '64.157.59'
$ws.Range("E16").Value = '  -3.47%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.115'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '3.142.01'
$ws.Range("E18").Value = '  -3.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.88'
$ws.Range("E19").Value = '  -7.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '477.53'
$ws.Range("E20").Value = '  -5.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.78'
$ws.Range("E21").Value = '  -4.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.707'
$ws.Range("E22").Value = '  -5.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.80'
$ws.Range("E23").Value = '  -3.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.64'
$ws.Range("E24").Value = '  -7.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.98'
$ws.Range("E25").Value = '  -4.73%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.89'
$ws.Range("E27").Value = '  -4.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.46'
$ws.Range("E28").Value = '  -7.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("E29").Value = '  -8.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.77'
$ws.Range("E30").Value = '  -3.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.114'
$ws.Range("E31").Value = '  -11.92%  '
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.75'
$ws.Range("E32").Value = '  -4.93%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.21'
$ws.Range("E34").Value = '  -6.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.12'
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.98'
$ws.Range("E36").Value = '  -7.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.47'
$ws.Range("E37").Value = '  -3.74%  '
$ws.Range("D38").Value = '0.0₃0738'
$ws.Range("E38").Value = '  -6.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '461.54'
$ws.Range("E39").Value = '  -6.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  -11.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0395'
$ws.Range("E41").Value = '  -6.18%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("E42").Value = '  -7.60%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.43'
$ws.Range("E43").Value = '  -4.30%  '
$ws.Range("D44").Value = '2.850.69'
$ws.Range("E44").Value = '  -4.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.266'
$ws.Range("E45").Value = '  -9.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.26'
$ws.Range("E46").Value = '  -9.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.47'
$ws.Range("E47").Value = '  -7.85%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("E49").Value = '  -6.32%  '
$ws.Range("E50").Value = '  -4.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.04'
$ws.Range("E51").Value = '  -1.76%  '
